$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Eva Gerecke): Highest Degree / Credits Beyond Degree / Date Available To Start
$ws.Range("L2").Value = "Master"

# "16" must stay text, not become a number - force the cell to Text format first
$ws.Range("M2").NumberFormat = "@"
$ws.Range("M2").Value = "16"

$ws.Range("N2").Value = "August 30, 2023"

# Row 3 (Setary Abbay): Highest Degree / Credits Beyond Degree / Date Available To Start
$ws.Range("L3").Value = "Bachelor"
$ws.Range("M3").Value = "None"
$ws.Range("N3").Value = "May 16, 2016"
